# Swap the order of "dnasr281@gmail.com" and "System" in the "Recorded By"
# column (G) wherever a session was recorded by both - i.e. change every
# "dnasr281@gmail.com, System" value to "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Output "Updated $changed 'Recorded By' cell(s) in column G."
